# Update generated "want to go" counts (column F) on the "展览" (Exhibition)
# and "全部类型" (All Types) sheets, matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$sheetExhibit = $wb.Worksheets.Item("展览")
$sheetExhibit.Range("F2").Value = 438
$sheetExhibit.Range("F3").Value = 5346
$sheetExhibit.Range("F6").Value = 49
$sheetExhibit.Range("F7").Value = 514

$sheetAll = $wb.Worksheets.Item("全部类型")
$sheetAll.Range("F2").Value = 438
$sheetAll.Range("F3").Value = 5346
$sheetAll.Range("F8").Value = 49
$sheetAll.Range("F9").Value = 514
